$d = $word.ActiveDocument
$d.Content.Find.Execute("shortest", $true, $false, $false, $false, $false,
                         $true, 1, $false, "minimum", 2)
